$wb = $excel.ActiveWorkbook

$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = -0.7217969836267406
$ws0.Range("C2").Value = -0.578672443733654
$ws0.Range("B3").Value = 0.340860819974193
$ws0.Range("C3").Value = 0.8711035324500809
$ws0.Range("B4").Value = 0.2337723897796488
$ws0.Range("C4").Value = -0.25820997514052

$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -1.350738518414303
$ws1.Range("C2").Value = -0.3001173235682472
$ws1.Range("B3").Value = 0.5520842911472739
$ws1.Range("C3").Value = 0.2739620574653207
$ws1.Range("B4").Value = 0.3296052199406059
$ws1.Range("C4").Value = 0.4263540533725379
